$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new 2-column table (Cineplex, CinemaID), with header row
$data = @(
    @("Cineplex", "CinemaID"),
    @("Downtown", "D01"),
    @("Downtown", "D02"),
    @("Downtown", "D03"),
    @("Causeway", "C01"),
    @("Causeway", "C02"),
    @("Causeway", "C03"),
    @("Tampines", "T01"),
    @("Tampines", "T02"),
    @("Tampines", "T03")
)

# Clear the old C and D columns (Cinema name / Type) entirely
$ws.Columns.Item(3).Clear()
$ws.Columns.Item(4).Clear()

# Write the new A:B table, including header row
for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $r + 1
    $ws.Cells.Item($row, 1).Value = $data[$r][0]
    $ws.Cells.Item($row, 2).Value = $data[$r][1]
}

# Set column B width to match new layout (closest reachable value to 10.140625)
$ws.Columns.Item(2).ColumnWidth = 9.25

# Update the selection to match the final saved state
$ws.Range("F8").Select()
